$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A19").Value = "'12/13/2025"
$ws.Range("A19").Style = "Normal"
$ws.Range("B19").Value = 12398.31
$ws.Range("C19").Value = 0.2015704448859508
$ws.Range("D19").Value = 0.7984295551140492
$ws.Range("E19").Value = -128.1
$ws.Range("F19").Value = -27.85
$ws.Range("G19").Value = -20617.58
$ws.Range("H19").Value = -67.56
$ws.Range("I19").Value = -403.72
$ws.Range("J19").Value = -13.91
